$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add the EPIC method column (H) with header and "x" marks on the same rows
# that already have an "x" in column G (xCell/MCPCounter/quantiseq style pattern).
$ws.Range("H1").Value = "EPIC"
$ws.Range("H1").Font.Bold = $true

$ws.Range("H9").Value = "x"
$ws.Range("H12").Value = "x"
$ws.Range("H13").Value = "x"
$ws.Range("H18").Value = "x"

# Make Tabelle1 the active sheet and select H17 on it, matching the saved view state.
$ws.Activate()
$ws.Range("H17").Select()
